$wb = $excel.ActiveWorkbook

# Reverse the "img" prefix/suffix naming convention: *img -> img*
$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# The active/selected tab moves from "holiday" (index 8) to the last
# sheet, now named "imge" (index 16).
$ws = $wb.Worksheets.Item("imge")
$ws.Activate()
